$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.772.25"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3
$ws.Range("D3").Value = "3.811.88"
$ws.Range("E3").Value = "  -1.40%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "709.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.70%  "

# Row 7
$ws.Range("D7").Value = "3.812.07"
$ws.Range("E7").Value = "  -1.36%  "

# Row 8
$ws.Range("E8").Value = "  +0.03%  "

# Row 9
$ws.Range("E9").Value = "  -0.63%  "

# Row 10
$ws.Range("E10").Value = "  -1.82%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.54"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.18%  "

# Row 12
$ws.Range("E12").Value = "  -0.66%  "

# Row 13
$ws.Range("E13").Value = "  -2.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.32%  "

# Row 15
$ws.Range("D15").Value = "4.454.09"
$ws.Range("E15").Value = "  -1.29%  "

# Row 16
$ws.Range("D16").Value = "3.829.41"
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").Value = "70.847.82"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18
$ws.Range("E18").Value = "  +0.18%  "

# Row 19
$ws.Range("E19").Value = "  -1.42%  "

# Row 20
$ws.Range("E20").Value = "  -2.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "496.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.62%  "

# Row 22
$ws.Range("E22").Value = "  -4.85%  "

# Row 23
$ws.Range("E23").Value = "  +0.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.25"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.06%  "

# Row 25
$ws.Range("E25").Value = "  -1.19%  "

# Row 26
$ws.Range("E26").Value = "  -1.95%  "

# Row 27
$ws.Range("E27").Value = "  -3.18%  "

# Row 28
$ws.Range("D28").Value = "3.962.24"
$ws.Range("E28").Value = "  -1.05%  "

# Row 29
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("E30").Value = "  -4.94%  "

# Row 31
$ws.Range("E31").Value = "  -3.28%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.23"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.33"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.12"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.171"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.85%  "

# Row 36
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.779.85"
$ws.Range("E36").Value = "  -0.96%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.11"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.13%  "

# Row 38
$ws.Range("E38").Value = "  -0.11%  "

# Row 39
$ws.Range("E39").Value = "  -2.93%  "

# Row 40
$ws.Range("E40").Value = "  +0.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.30"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.16%  "

# Row 42
$ws.Range("E42").Value = "  -2.08%  "

# Row 43
$ws.Range("E43").Value = "  -4.89%  "

# Row 44
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("E45").Value = "  +0.22%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000320"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "166.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.81"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.41%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "423.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.02%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.62"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("E51").Value = "  -3.35%  "
